$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D target cells to be treated as Text so that numeric-looking
# strings (e.g. "1.00", "0.0000235") keep their exact original formatting
# instead of being reinterpreted as numbers by Excel.
$dCells = @("D2","D3","D5","D6","D9","D10","D13","D14","D16","D17","D18","D19","D20","D22","D24","D25","D26","D27","D28","D30","D33","D35","D36","D38","D39","D40","D41","D44","D45","D46","D48","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply cell value updates
$ws.Range("D2").Value = "61.795.09"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").Value = "2.908.21"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "586.13"
$ws.Range("E5").Value = "  -1.53%  "
$ws.Range("D6").Value = "146.43"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "2.907.84"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").Value = "6.84"
$ws.Range("E10").Value = "  -6.79%  "
$ws.Range("E11").Value = "  +5.19%  "
$ws.Range("E12").Value = "  -3.18%  "
$ws.Range("D13").Value = "0.0000235"
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").Value = "32.75"
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("E15").Value = "  -0.98%  "
$ws.Range("D16").Value = "3.389.52"
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("D17").Value = "61.816.18"
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("D18").Value = "6.65"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").Value = "2.908.22"
$ws.Range("E19").Value = "  -2.09%  "
$ws.Range("D20").Value = "435.44"
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").Value = "0.659"
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("D24").Value = "80.54"
$ws.Range("E24").Value = "  -1.67%  "
$ws.Range("D25").Value = "11.95"
$ws.Range("D26").Value = "10.19"
$ws.Range("E26").Value = "  -9.02%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").Value = "2.07"
$ws.Range("E28").Value = "  -2.90%  "
$ws.Range("E29").Value = "  +21.58%  "
$ws.Range("D30").Value = "7.13"
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("E31").Value = "  -1.68%  "
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("D33").Value = "0.109"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "25.85"
$ws.Range("E35").Value = "  -2.99%  "
$ws.Range("D36").Value = "0.971"
$ws.Range("E36").Value = "  -2.20%  "
$ws.Range("E37").Value = "  -2.55%  "
$ws.Range("D38").Value = "3.02"
$ws.Range("E38").Value = "  +3.49%  "
$ws.Range("D39").Value = "49.19"
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("D40").Value = "1.99"
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("D41").Value = "8.37"
$ws.Range("E41").Value = "  -2.21%  "
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("E43").Value = "  -3.42%  "
$ws.Range("D44").Value = "38.52"
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("D45").Value = "2.699.19"
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("D46").Value = "134.01"
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("E47").Value = "  -1.71%  "
$ws.Range("D48").Value = "342.03"
$ws.Range("E48").Value = "  -6.05%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  -1.78%  "
$ws.Range("D51").Value = "22.18"
$ws.Range("E51").Value = "  -2.93%  "
